# Update the "想去人数" (interested-count) figures in the 展览 and 全部类型
# sheets to reflect the newly generated numbers.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row => new F value for sheet "展览"
$sheet1Updates = @{
    2  = 329
    3  = 280
    4  = 1214
    10 = 3417
    17 = 76
    18 = 705
    24 = 2513
    25 = 5028
    29 = 1288
    31 = 2212
    32 = 570
    38 = 456
    39 = 781
    40 = 28
}

# Row => new F value for sheet "全部类型"
$sheet4Updates = @{
    2  = 329
    3  = 280
    4  = 1214
    10 = 3417
    18 = 76
    19 = 705
    25 = 2513
    26 = 5028
    30 = 1288
    32 = 2212
    33 = 570
    39 = 456
    40 = 781
    41 = 28
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
